# comment from python script
#
# Ninja sheet (4th tab): a new "T3" (column D) score was recorded for every
# student. Column P ("Resultados") already holds a shared
# IFERROR(SUM(..)/COUNT(..)*100,0) formula, so it recalculates on its own
# once column D is populated. Row 23 (a stray formatted-but-empty D23 cell
# below the table) is removed as part of the same cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ninja")

# Row -> T3 (column D) value, in row order (rows 2-21).
$t3Scores = @(0, 1, 1, 1, 1, 1, 1, 1, 0, 0, 1, 0, 1, 1, 1, 1, 1, 1, 1, 1)

$startRow = 2
for ($i = 0; $i -lt $t3Scores.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 4).Value = $t3Scores[$i]
}

# Remove the now-obsolete formatted cell/row sitting below the data table.
$ws.Rows.Item(23).Delete()

$ws.Range("E2").Select()
